$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the target sentence: "Es una máquina virtual de sistema."
# ------------------------------------------------------------------
$target = $d.Content
$found = $target.Find.Execute("Es una máquina virtual de sistema.", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $found) {
    throw "Could not locate the target sentence."
}

$sentenceStart = $target.Start
$sentenceEnd = $target.End
$textColor = $target.Font.Color

# ------------------------------------------------------------------
# Update paragraph indentation: left=1440/firstLine=720 -> left=2160 (no firstLine)
# ------------------------------------------------------------------
$target.ParagraphFormat.LeftIndent = 108
$target.ParagraphFormat.FirstLineIndent = 0

# ------------------------------------------------------------------
# Split "Es una máquina virtual de sistema." into three runs:
#   1) "Es una máquina virtual de sistema"
#   2) ", es decir, permite realizar instalaciones de sistemas operativos en"
#   3) " diferentes instancias virtuales."
# ------------------------------------------------------------------

# Drop the trailing period from the first run.
$periodRange = $d.Range($sentenceEnd - 1, $sentenceEnd)
$periodRange.Text = ""

$run1End = $sentenceEnd - 1

# Insert the second run's text right after the first run.
$secondText = ", es decir, permite realizar instalaciones de sistemas operativos en"
$insertPoint2 = $d.Range($run1End, $run1End)
$insertPoint2.InsertAfter($secondText)
$run2Start = $run1End
$run2End = $run2Start + $secondText.Length
$run2Range = $d.Range($run2Start, $run2End)
$run2Range.Font.Color = $textColor
# Force this insertion to stay a distinct run rather than being
# re-coalesced into the preceding (identically formatted) run.
$run2Range.Font.Bold = $true
$run2Range.Font.Bold = $false

# Insert the third run's text right after the second run.
$thirdText = " diferentes instancias virtuales."
$insertPoint3 = $d.Range($run2End, $run2End)
$insertPoint3.InsertAfter($thirdText)
$run3Start = $run2End
$run3End = $run3Start + $thirdText.Length
$run3Range = $d.Range($run3Start, $run3End)
$run3Range.Font.Color = $textColor
$run3Range.Font.Bold = $true
$run3Range.Font.Bold = $false

$finalRange = $d.Range($sentenceStart, $run3End)
Write-Host "Final paragraph text: [$($finalRange.Text)]"
